# "Modifique detalles en municipios"
#
# - Row 12 used to hold a single left-over formatted (but empty) cell at
#   L12. Replace it with a filled-in helper/check row that numbers every
#   column from A12 (=1) through AM12 (=39), one more than the previous
#   column - built the same way a user would in Excel: type the first
#   couple of formulas by hand and then fill the rest of the row to the
#   right.
# - L12 no longer needs its old one-off underline style, so clear its
#   direct formatting back to Normal.
# - Update the sheet view: the sheet was scrolled/selected near the end of
#   the "Clave Municipio" column; move the selection over to the (now
#   last/whole) AM column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray formatting that used to sit (empty) in L12.
$ws.Range("L12").ClearFormats()

# Seed the first three cells by hand ...
$ws.Range("A12").Value = 1
$ws.Range("B12").Formula = "=A12+1"
$ws.Range("C12").Formula = "=B12+1"

# ... then fill the same "+1" pattern across the rest of the row (D:AM),
# exactly the way Excel's own Fill Right would continue the series.
$ws.Range("D12:AM12").Formula = "=C12+1"

# Re-point the view: select the full Clave Municipio column (AM) and bring
# it into frame.
$ws.Range("AM1:AM1048576").Select()

$wb.Save()
